$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.82499266666667
$ws.Range("H2").Value = 56.474978
$ws.Range("I2").Value = 0.06886869772378311
$ws.Range("J2").Value = 0.0688686977237831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.815493333333333
$ws.Range("N2").Value = 5.44648
$ws.Range("O2").Value = 0.02449420924905278
$ws.Range("P2").Value = 0.02449420924905277
$ws.Range("Q2").Value = 34.17664868638222
$ws.Range("R2").Value = 307.58983817744
$ws.Range("S2").Value = 0.001686884292756108
$ws.Range("T2").Value = 0.001686884292756108

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.82499266666667
$ws.Range("H3").Value = 56.474978
$ws.Range("I3").Value = 0.06886869772378311
$ws.Range("J3").Value = 0.0688686977237831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.684019666666667
$ws.Range("N3").Value = 17.052059
$ws.Range("O3").Value = 0.07668745708663094
$ws.Range("P3").Value = 0.07668745708663093
$ws.Range("Q3").Value = 107.0016285421891
$ws.Range("R3").Value = 963.014656879702
$ws.Range("S3").Value = 0.005281365301304775
$ws.Range("T3").Value = 0.005281365301304773

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.82499266666667
$ws.Range("H4").Value = 56.474978
$ws.Range("I4").Value = 0.06886869772378311
$ws.Range("J4").Value = 0.0688686977237831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 61.37607633333332
$ws.Range("N4").Value = 184.128229
$ws.Range("O4").Value = 0.8280715929891429
$ws.Range("P4").Value = 0.8280715929891428
$ws.Range("Q4").Value = 1155.404186883773
$ws.Range("R4").Value = 10398.63768195396
$ws.Range("S4").Value = 0.05702821223122084
$ws.Range("T4").Value = 0.05702821223122082

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.82499266666667
$ws.Range("H5").Value = 56.474978
$ws.Range("I5").Value = 0.06886869772378311
$ws.Range("J5").Value = 0.0688686977237831
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.243698
$ws.Range("N5").Value = 15.731094
$ws.Range("O5").Value = 0.07074674067517345
$ws.Range("P5").Value = 0.07074674067517343
$ws.Range("Q5").Value = 98.71257639621467
$ws.Range("R5").Value = 888.413187565932
$ws.Range("S5").Value = 0.004872235898501391
$ws.Range("T5").Value = 0.00487223589850139

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 121.8208923333333
$ws.Range("H6").Value = 365.462677
$ws.Range("I6").Value = 0.4456653109566078
$ws.Range("J6").Value = 0.4456653109566078
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.815493333333333
$ws.Range("N6").Value = 5.44648
$ws.Range("O6").Value = 0.02449420924905278
$ws.Range("P6").Value = 0.02449420924905277
$ws.Range("Q6").Value = 221.1650178918845
$ws.Range("R6").Value = 1990.48516102696
$ws.Range("S6").Value = 0.01091621938161533
$ws.Range("T6").Value = 0.01091621938161532

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 121.8208923333333
$ws.Range("H7").Value = 365.462677
$ws.Range("I7").Value = 0.4456653109566078
$ws.Range("J7").Value = 0.4456653109566078
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.684019666666667
$ws.Range("N7").Value = 17.052059
$ws.Range("O7").Value = 0.07668745708663094
$ws.Range("P7").Value = 0.07668745708663093
$ws.Range("Q7").Value = 692.4323478335493
$ws.Range("R7").Value = 6231.891130501944
$ws.Range("S7").Value = 0.0341769394089849
$ws.Range("T7").Value = 0.03417693940898488

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 121.8208923333333
$ws.Range("H8").Value = 365.462677
$ws.Range("I8").Value = 0.4456653109566078
$ws.Range("J8").Value = 0.4456653109566078
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 61.37607633333332
$ws.Range("N8").Value = 184.128229
$ws.Range("O8").Value = 0.8280715929891429
$ws.Range("P8").Value = 0.8280715929891428
$ws.Range("Q8").Value = 7476.888386845448
$ws.Range("R8").Value = 67291.99548160903
$ws.Range("S8").Value = 0.36904278398384
$ws.Range("T8").Value = 0.3690427839838399

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 121.8208923333333
$ws.Range("H9").Value = 365.462677
$ws.Range("I9").Value = 0.4456653109566078
$ws.Range("J9").Value = 0.4456653109566078
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.243698
$ws.Range("N9").Value = 15.731094
$ws.Range("O9").Value = 0.07074674067517345
$ws.Range("P9").Value = 0.07074674067517345
$ws.Range("Q9").Value = 638.7919694865154
$ws.Range("R9").Value = 5749.127725378639
$ws.Range("S9").Value = 0.03152936818216767
$ws.Range("T9").Value = 0.03152936818216766

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 87.673585
$ws.Range("H10").Value = 263.020755
$ws.Range("I10").Value = 0.3207419907481189
$ws.Range("J10").Value = 0.3207419907481188
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.815493333333333
$ws.Range("N10").Value = 5.44648
$ws.Range("O10").Value = 0.02449420924905278
$ws.Range("P10").Value = 0.02449420924905277
$ws.Range("Q10").Value = 159.1708090769334
$ws.Range("R10").Value = 1432.5372816924
$ws.Range("S10").Value = 0.007856321436342174
$ws.Range("T10").Value = 0.007856321436342172

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 87.673585
$ws.Range("H11").Value = 263.020755
$ws.Range("I11").Value = 0.3207419907481189
$ws.Range("J11").Value = 0.3207419907481188
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.684019666666667
$ws.Range("N11").Value = 17.052059
$ws.Range("O11").Value = 0.07668745708663094
$ws.Range("P11").Value = 0.07668745708663093
$ws.Range("Q11").Value = 498.3383813871717
$ws.Range("R11").Value = 4485.045432484545
$ws.Range("S11").Value = 0.02459688765137694
$ws.Range("T11").Value = 0.02459688765137694

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 87.673585
$ws.Range("H12").Value = 263.020755
$ws.Range("I12").Value = 0.3207419907481189
$ws.Range("J12").Value = 0.3207419907481188
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 61.37607633333332
$ws.Range("N12").Value = 184.128229
$ws.Range("O12").Value = 0.8280715929891429
$ws.Range("P12").Value = 0.8280715929891428
$ws.Range("Q12").Value = 5381.060645376988
$ws.Range("R12").Value = 48429.54580839289
$ws.Range("S12").Value = 0.2655973312173037
$ws.Range("T12").Value = 0.2655973312173037

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 87.673585
$ws.Range("H13").Value = 263.020755
$ws.Range("I13").Value = 0.3207419907481189
$ws.Range("J13").Value = 0.3207419907481188
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.243698
$ws.Range("N13").Value = 15.731094
$ws.Range("O13").Value = 0.07074674067517345
$ws.Range("P13").Value = 0.07074674067517343
$ws.Range("Q13").Value = 459.73380231733
$ws.Range("R13").Value = 4137.60422085597
$ws.Range("S13").Value = 0.02269145044309605
$ws.Range("T13").Value = 0.02269145044309604

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 45.02666966666666
$ws.Range("H14").Value = 135.080009
$ws.Range("I14").Value = 0.1647240005714903
$ws.Range("J14").Value = 0.1647240005714903
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.815493333333333
$ws.Range("N14").Value = 5.44648
$ws.Range("O14").Value = 0.02449420924905278
$ws.Range("P14").Value = 0.02449420924905277
$ws.Range("Q14").Value = 81.74561860203555
$ws.Range("R14").Value = 735.71056741832
$ws.Range("S14").Value = 0.004034784138339172
$ws.Range("T14").Value = 0.004034784138339171

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 45.02666966666666
$ws.Range("H15").Value = 135.080009
$ws.Range("I15").Value = 0.1647240005714903
$ws.Range("J15").Value = 0.1647240005714903
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.684019666666667
$ws.Range("N15").Value = 17.052059
$ws.Range("O15").Value = 0.07668745708663094
$ws.Range("P15").Value = 0.07668745708663093
$ws.Range("Q15").Value = 255.9324759098368
$ws.Range("R15").Value = 2303.392283188531
$ws.Range("S15").Value = 0.01263226472496433
$ws.Range("T15").Value = 0.01263226472496433

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 45.02666966666666
$ws.Range("H16").Value = 135.080009
$ws.Range("I16").Value = 0.1647240005714903
$ws.Range("J16").Value = 0.1647240005714903
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 61.37607633333332
$ws.Range("N16").Value = 184.128229
$ws.Range("O16").Value = 0.8280715929891429
$ws.Range("P16").Value = 0.8280715929891428
$ws.Range("Q16").Value = 2763.560314497117
$ws.Range("R16").Value = 24872.04283047406
$ws.Range("S16").Value = 0.1364032655567785
$ws.Range("T16").Value = 0.1364032655567784

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 45.02666966666666
$ws.Range("H17").Value = 135.080009
$ws.Range("I17").Value = 0.1647240005714903
$ws.Range("J17").Value = 0.1647240005714903
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.243698
$ws.Range("N17").Value = 15.731094
$ws.Range("O17").Value = 0.07074674067517345
$ws.Range("P17").Value = 0.07074674067517343
$ws.Range("Q17").Value = 236.1062576777607
$ws.Range("R17").Value = 2124.956319099846
$ws.Range("S17").Value = 0.01165368615140835
$ws.Range("T17").Value = 0.01165368615140834
